$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oni")

# Extend the yearly table with a new 2023 column (K), mirroring the
# formatting already used for the 2022 column (J).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 757.2
$ws.Range("K5").Value = 940.6
$ws.Range("K6").Value = 687.9
